$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 and 39: coin name/link swap (order changed) plus price/volume updates
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"

# Update Price (D) and Volume(1h) (E) columns for each row

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.054.71"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.765.87"
$ws.Range("E3").Value = "  -1.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.37"
$ws.Range("E5").Value = "  -1.22%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.08%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3912"
$ws.Range("E7").Value = "  +1.87%  "

# Row 8
$ws.Range("E8").Value = "  -1.22%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.10"
$ws.Range("E9").Value = "  -4.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.125"
$ws.Range("E10").Value = "  -2.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07244"
$ws.Range("E11").Value = "  -2.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.52"
$ws.Range("E12").Value = "  -3.67%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  +0.15%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.158"
$ws.Range("E14").Value = "  -4.77%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.132"
$ws.Range("E15").Value = "  -3.29%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.760.81"
$ws.Range("E16").Value = "  -1.20%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001060"
$ws.Range("E17").Value = "  -1.76%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06622"
$ws.Range("E18").Value = "  -1.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.22"
$ws.Range("E19").Value = "  -2.35%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9990"
$ws.Range("E20").Value = "  -0.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.98"
$ws.Range("E21").Value = "  -2.90%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.233"
$ws.Range("E22").Value = "  -3.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.050.26"
$ws.Range("E23").Value = "  -0.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.68"
$ws.Range("E24").Value = "  -3.51%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.380"
$ws.Range("E25").Value = "  +0.93%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.14"
$ws.Range("E26").Value = "  +0.60%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.00"
$ws.Range("E27").Value = "  -3.77%  "

# Row 28
$ws.Range("E28").Value = "  -3.43%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.963.14"
$ws.Range("E29").Value = "  -0.58%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.290"
$ws.Range("E30").Value = "  -9.64%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.72"
$ws.Range("E31").Value = "  -4.35%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.082"
$ws.Range("E32").Value = "  +1.67%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.829"
$ws.Range("E33").Value = "  -4.75%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08746"
$ws.Range("E34").Value = "  -2.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.11"
$ws.Range("E35").Value = "  -5.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06218"
$ws.Range("E36").Value = "  -2.96%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02294"
$ws.Range("E37").Value = "  -5.52%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6520"
$ws.Range("E38").Value = "  -5.28%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.136"
$ws.Range("E39").Value = "  -4.49%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2120"
$ws.Range("E40").Value = "  -2.40%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.498"
$ws.Range("E41").Value = "  -0.33%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.204"
$ws.Range("E42").Value = "  -3.47%  "

# Row 43
$ws.Range("E43").Value = "  -4.80%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9986"
$ws.Range("E44").Value = "  -0.22%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.91"
$ws.Range("E45").Value = "  -1.69%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.830"
$ws.Range("E46").Value = "  -1.21%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6018"
$ws.Range("E47").Value = "  -4.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.03"
$ws.Range("E48").Value = "  -4.56%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.002"
$ws.Range("E49").Value = "  -3.99%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.157"
$ws.Range("E50").Value = "  -4.30%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07024"
$ws.Range("E51").Value = "  -6.19%  "
